$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.885.79"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.767.69"
$cell.Style = "Normal"

$ws.Range("E3").Value = "  +7.24%  "
$ws.Range("E4").Value = "  -0.12%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "421.39"
$cell.Style = "Normal"

$ws.Range("E5").Value = "  +0.64%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "132.47"
$cell.Style = "Normal"

$ws.Range("E6").Value = "  -0.05%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.750.55"
$cell.Style = "Normal"

$ws.Range("E7").Value = "  +6.94%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.653"
$cell.Style = "Normal"

$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +0.08%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.776"
$cell.Style = "Normal"

$ws.Range("E10").Value = "  -0.40%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.188"
$cell.Style = "Normal"

$ws.Range("E11").Value = "  +15.22%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.0000435"
$cell.Style = "Normal"

$ws.Range("E12").Value = "  +63.51%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "43.01"
$cell.Style = "Normal"

$ws.Range("E13").Value = "  -0.93%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "10.40"
$cell.Style = "Normal"

$ws.Range("E14").Value = "  +3.81%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.358.99"
$cell.Style = "Normal"

$ws.Range("E15").Value = "  +7.24%  "
$ws.Range("E16").Value = "  -0.51%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.776.79"
$cell.Style = "Normal"

$ws.Range("E17").Value = "  +8.25%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "20.64"
$cell.Style = "Normal"

$ws.Range("E18").Value = "  +0.38%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "13.07"
$cell.Style = "Normal"

$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("E20").Value = "  +3.18%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "67.889.43"
$cell.Style = "Normal"

$ws.Range("E21").Value = "  +3.82%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "452.24"
$cell.Style = "Normal"

$ws.Range("E22").Value = "  -0.77%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "15.97"
$cell.Style = "Normal"

$ws.Range("E23").Value = "  +19.84%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "90.93"
$cell.Style = "Normal"

$ws.Range("E24").Value = "  +0.69%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "3.11"
$cell.Style = "Normal"

$ws.Range("E25").Value = "  -4.34%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "38.53"
$cell.Style = "Normal"

$ws.Range("E26").Value = "  +12.44%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "3.36"
$cell.Style = "Normal"

$ws.Range("E27").Value = "  -1.37%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.14"
$cell.Style = "Normal"

$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("E29").Value = "  +5.40%  "
$ws.Range("E30").Value = "  +6.15%  "
$ws.Range("E31").Value = "  +0.98%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.70"
$cell.Style = "Normal"

$ws.Range("E32").Value = "  -1.53%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "7.20"
$cell.Style = "Normal"

$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("E34").Value = "  +1.76%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "42.20"
$cell.Style = "Normal"

$ws.Range("E35").Value = "  +6.14%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "57.71"
$cell.Style = "Normal"

$ws.Range("E36").Value = "  +0.33%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"

$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -2.31%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0760"
$cell.Style = "Normal"

$ws.Range("E39").Value = "  +3.14%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.03"
$cell.Style = "Normal"

$ws.Range("E40").Value = "  +30.49%  "
$ws.Range("E41").Value = "  +0.38%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"

$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +29.05%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.39"
$cell.Style = "Normal"

$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  +6.35%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "147.73"
$cell.Style = "Normal"

$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.16"
$cell.Style = "Normal"

$ws.Range("E47").Value = "  +22.97%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.Style = "Normal"

$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "4.35"
$cell.Style = "Normal"

$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.65"
$cell.Style = "Normal"

$ws.Range("E50").Value = "  -3.87%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.310"
$cell.Style = "Normal"

$ws.Range("E51").Value = "  -0.85%  "
